$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update existing row 24: fill in Description ("Testing Required") and
# change Runmode from "Y" to "Done"
$ws.Range("B24").Value = "Testing Required"
$ws.Range("C24").Value = "Done"

# Add new row 25 for the ImportMojio test case
$ws.Range("A25").Value = "ImportMojio"
$ws.Range("B25").Value = "Testing Required"
$ws.Range("C25").Value = "Y"

# Match the style used on the rest of the data rows
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null

# Update the active cell selection to A25, matching the new state
$ws.Range("A25").Select() | Out-Null
